$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper behavior: for numeric-looking text values (Price / Volume columns),
# force the cell to Text format before assignment so Excel does not auto-convert
# the string into a number, then restore the General number format to match the
# original workbook formatting (values remain stored as text/inlineStr).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.607.03'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("E2").NumberFormat = "General"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.268.61'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E3").NumberFormat = "General"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E4").NumberFormat = "General"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.69'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.80%  '
$ws.Range("E5").NumberFormat = "General"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.85'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("E6").NumberFormat = "General"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E7").NumberFormat = "General"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.265.41'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("E8").NumberFormat = "General"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.569'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("E9").NumberFormat = "General"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E10").NumberFormat = "General"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.570'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("E11").NumberFormat = "General"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.06'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("E12").NumberFormat = "General"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000262'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.50%  '
$ws.Range("E13").NumberFormat = "General"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.795.30'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("E14").NumberFormat = "General"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.37'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("E15").NumberFormat = "General"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '611.60'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.45%  '
$ws.Range("E16").NumberFormat = "General"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.529.30'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("E17").NumberFormat = "General"

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E18").NumberFormat = "General"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.77'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("E19").NumberFormat = "General"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.294.54'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E20").NumberFormat = "General"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.87'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.95%  '
$ws.Range("E21").NumberFormat = "General"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.886'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.05%  '
$ws.Range("E22").NumberFormat = "General"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.96'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E23").NumberFormat = "General"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.01'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.97%  '
$ws.Range("E24").NumberFormat = "General"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.94'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E25").NumberFormat = "General"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E26").NumberFormat = "General"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E27").NumberFormat = "General"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.45'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E28").NumberFormat = "General"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.73'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("E29").NumberFormat = "General"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("E30").NumberFormat = "General"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.42'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E31").NumberFormat = "General"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.72'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.33%  '
$ws.Range("E32").NumberFormat = "General"

# Row 33
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.81'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("E33").NumberFormat = "General"

# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '541.34'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E34").NumberFormat = "General"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.779.50'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("E35").NumberFormat = "General"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.103'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("E36").NumberFormat = "General"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.997'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("E37").NumberFormat = "General"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.99'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.70%  '
$ws.Range("E38").NumberFormat = "General"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.127'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E39").NumberFormat = "General"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '32.42'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("E40").NumberFormat = "General"

# Row 41
$ws.Range("B41").Value = 'ApeXProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.39'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.86%  '
$ws.Range("E41").NumberFormat = "General"

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.12'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.89%  '
$ws.Range("E42").NumberFormat = "General"

# Row 43
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0675'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.30%  '
$ws.Range("E43").NumberFormat = "General"

# Row 44
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.37%  '
$ws.Range("E44").NumberFormat = "General"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.329'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("E45").NumberFormat = "General"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0404'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("E46").NumberFormat = "General"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.99'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.11%  '
$ws.Range("E47").NumberFormat = "General"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("E48").NumberFormat = "General"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("E49").NumberFormat = "General"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.49'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.01%  '
$ws.Range("E50").NumberFormat = "General"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '127.84'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.65%  '
$ws.Range("E51").NumberFormat = "General"
